$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.119.44'
$ws.Range("E2").Value = '  -2.03%  '

$ws.Range("D3").Value = '2.638.19'
$ws.Range("E3").Value = '  -3.16%  '

$ws.Range("E4").Value = '  -0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.20'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.21%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.10%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.543'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.89%  '

$ws.Range("D9").Value = '2.637.84'
$ws.Range("E9").Value = '  -3.15%  '

$ws.Range("E10").Value = '  -1.28%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.160'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.48%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.363'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.95%  '

$ws.Range("E13").Value = '  -1.14%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.28%  '

$ws.Range("D15").Value = '3.117.37'
$ws.Range("E15").Value = '  -3.06%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000182'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.78%  '

$ws.Range("D17").Value = '67.188.14'
$ws.Range("E17").Value = '  -1.96%  '

$ws.Range("D18").Value = '2.635.11'
$ws.Range("E18").Value = '  -0.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.49%  '

$ws.Range("E20").Value = '  +6.48%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '361.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.36'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.46%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.77'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.70%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +8.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.21%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '70.73'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.39%  '

$ws.Range("D28").Value = '2.770.57'
$ws.Range("E28").Value = '  -3.34%  '

$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.73%  '

$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000102'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.68%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '554.70'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.87%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.97'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.35%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.38'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.91'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.23%  '

$ws.Range("E35").Value = '  +1.85%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.53'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '157.75'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.26'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.370'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.52%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.27'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.60%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.88'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.02%  '

$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.51'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.76%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.97%  '

$ws.Range("B47").Value = 'ARBITRUM'
$ws.Range("C47").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.592'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.23%  '

$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").Value = '0.0₆0297'
$ws.Range("E48").Value = '  -4.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '153.07'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.85'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.73'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.15%  '
